# LOQ4249.docx edit:
#   1. "Programa" paragraph: split the single run into three <w:t> runs of
#      text joined by manual line breaks (<w:br/>), all inside one <w:r>.
#   2/3. "Avaliação" bullet paragraph: split the "Método:" and "Critério:"
#      value runs each into two <w:t> runs joined by a manual line break,
#      keeping the bold label runs untouched.
#   4. "Bibliografia" paragraph: split the single run into six <w:t> runs
#      (one per reference) inside one <w:r>, each pair joined by two manual
#      line breaks (<w:br/><w:br/>).
#
# Range.InsertXML() on a sub-run Range (one that doesn't start/end exactly
# at a paragraph boundary) misplaces the inserted content in this host, so
# every edit below rebuilds the *whole* paragraph (its <w:pPr> plus every
# <w:r>) and feeds it through InsertXML on that paragraph's full Range,
# which reliably replaces the paragraph's content in place.
#
# Also: calling InsertXML on the very last body paragraph (immediately
# before <w:sectPr>) leaves behind a stray empty trailing paragraph; the
# helper below detects and removes that artifact when it happens.

$d = $word.ActiveDocument

function Get-ParagraphContaining($doc, [string]$needle) {
    $match = $null
    foreach ($para in $doc.Paragraphs) {
        if ($para.Range.Text.Contains($needle)) {
            $match = $para
        }
    }
    return $match
}

function Set-ParagraphXml($doc, $paragraph, [string]$pPrXml, [string]$runsXml) {
    $pre = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p>'
    $post = '</w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $xml = $pre + $pPrXml + $runsXml + $post

    $countBefore = $doc.Paragraphs.Count
    $paragraph.Range.InsertXML($xml)
    $countAfter = $doc.Paragraphs.Count

    if ($countAfter -gt $countBefore) {
        # Editing the final body paragraph stranded an empty paragraph
        # after it (right before sectPr) -- remove the extra paragraph mark.
        $newEnd = $paragraph.Range.End
        $docEnd = $doc.Content.End
        $stray = $doc.Range($newEnd - 1, $docEnd)
        $stray.Delete()
    }
}

# --- Change 1: "Programa" paragraph ---------------------------------------
$p1 = Get-ParagraphContaining $d "Os mentores realizam apresentações sobre"
$runs1 = '<w:r>' +
         '<w:t xml:space="preserve">Proposição pelos alunos de startup de base tecnológica </w:t>' +
         '<w:br/>' +
         '<w:t xml:space="preserve">A proposta de startup é acompanhada por uma equipe de mentores, coordenada pelo professor da disciplina. </w:t>' +
         '<w:br/>' +
         '<w:t>Os mentores realizam apresentações sobre: inovação em produtos e serviços; necessidades e comportamento dos usuários; técnicas de ideação; definição de mercados, rotas tecnológicas e noções de propriedade intelectual; inovação aberta, capital de risco e técnicas de pitch</w:t>' +
         '</w:r>'
Set-ParagraphXml $d $p1 '' $runs1

# --- Changes 2 & 3: "Avaliação" bullet paragraph (Método / Critério) ------
$p2 = Get-ParagraphContaining $d "Atividades docentes: Mentoria"
$pPr2 = '<w:pPr><w:pStyle w:val="ListBullet"/></w:pPr>'
$runs2 = '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Método: </w:t></w:r>' +
         '<w:r><w:t>Atividades docentes: Mentoria, palestras e seminários.</w:t><w:br/><w:t>Atividades discentes: Elaboração de projeto utilizando laboratórios e instalações da USP.</w:t><w:br/></w:r>' +
         '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Critério: </w:t></w:r>' +
         '<w:r><w:t>Avaliação pela equipe de mentores, considerando critérios, tais como: qualidade técnica da proposta, grau de inovação, viabilidade técnica, dentre outros.</w:t><w:br/><w:t>Nota de projeto maior ou igual a 5,0 (cinco).</w:t><w:br/></w:r>' +
         '<w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Norma de recuperação: </w:t></w:r>' +
         '<w:r><w:t>Não há recuperação.</w:t></w:r>'
Set-ParagraphXml $d $p2 $pPr2 $runs2

# --- Change 4: "Bibliografia" paragraph ------------------------------------
$p3 = Get-ParagraphContaining $d "BROWN, T. Design thinking"
$runs3 = '<w:r>' +
         '<w:t>BROWN, T. Design thinking. Rio de Janeiro: Campus, 2010.</w:t><w:br/><w:br/>' +
         '<w:t>INPI. Instituto nacional de propriedade industrial. Disponível em: http://www.inpi.gov.br/. Consultado em: junho de 2015.</w:t><w:br/><w:br/>' +
         '<w:t>KUMAR, V. 101 Design Methods: A Structured Approach for Driving Innovation in Your Organization. New Jersey: John Willey and Sons, 2013.</w:t><w:br/><w:br/>' +
         '<w:t>MALHOTRA, N.K. Pesquisa de marketing: uma orientação aplicada. Porto Alegre: Bookman, 2006.</w:t><w:br/><w:br/>' +
         '<w:t>ROMEIRO FILHO et al. Projeto do produto. Rio de Janeiro: Campus, 2010.</w:t><w:br/><w:br/>' +
         '<w:t>ROZENFELD, Henrique. Gestão de desenvolvimento de produtos: uma referência para a melhoria do processo. Saraiva, 2006.</w:t>' +
         '</w:r>'
Set-ParagraphXml $d $p3 '' $runs3

Write-Output "edits applied"
